# Applies updated evaluation metric values (Value (average) and Value (std))
# across the RF, LGBM, XGB, KNN, and SVM worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("RF")
$ws.Range("C2").Value = 0.7137416098462077
$ws.Range("D2").Value = 0.05019111003300456
$ws.Range("C3").Value = 0.898743455497382
$ws.Range("D3").Value = 0.01606285202188616
$ws.Range("C4").Value = 0.8620498369646896
$ws.Range("D4").Value = 0.07326234427107983
$ws.Range("C5").Value = 0.5051207346574993
$ws.Range("D5").Value = 0.08520665744901536
$ws.Range("C6").Value = 0.9821
$ws.Range("D6").Value = 0.01155750485041416
$ws.Range("C7").Value = 0.6322200896668806
$ws.Range("D7").Value = 0.07228887897531978
$ws.Range("C8").Value = 0.8872344691826819
$ws.Range("D8").Value = 0.0196144968532882
$ws.Range("C9").Value = 0.7867248805998628
$ws.Range("D9").Value = 0.04030718596036099
$ws.Range("C10").Value = 0.7436111899269156
$ws.Range("D10").Value = 0.04195348795330237
$ws.Range("C11").Value = 0.6089300837940356
$ws.Range("D11").Value = 0.07088892424805661
$ws.Range("C12").Value = 0.903814
$ws.Range("D12").Value = 0.01473492325277963
$ws.Range("C13").Value = 0.7436111899269156
$ws.Range("D13").Value = 0.04195348795330234

$ws = $wb.Worksheets.Item("LGBM")
$ws.Range("C2").Value = 0.7237708077981441
$ws.Range("D2").Value = 0.04431002561367556
$ws.Range("C3").Value = 0.9
$ws.Range("D3").Value = 0.01805163558137198
$ws.Range("C4").Value = 0.839621112007205
$ws.Range("D4").Value = 0.08428257878115385
$ws.Range("C5").Value = 0.535662910618793
$ws.Range("D5").Value = 0.09164656429559094
$ws.Range("C6").Value = 0.977162
$ws.Range("D6").Value = 0.01410768917567832
$ws.Range("C7").Value = 0.6485324070412116
$ws.Range("D7").Value = 0.07538205089187128
$ws.Range("C8").Value = 0.8904213490170362
$ws.Range("D8").Value = 0.02112885433552696
$ws.Range("C9").Value = 0.7950841674140379
$ws.Range("D9").Value = 0.0424146750417171
$ws.Range("C10").Value = 0.7564082282907603
$ws.Range("D10").Value = 0.04518290398574542
$ws.Range("C11").Value = 0.6173994915816988
$ws.Range("D11").Value = 0.07711209846562808
$ws.Range("C12").Value = 0.9088540000000001
$ws.Range("D12").Value = 0.01613892758746595
$ws.Range("C13").Value = 0.7564082282907603
$ws.Range("D13").Value = 0.04518290398574541

$ws = $wb.Worksheets.Item("XGB")
$ws.Range("C2").Value = 0.7436888401182444
$ws.Range("D2").Value = 0.04257422999643389
$ws.Range("C3").Value = 0.9012565445026177
$ws.Range("D3").Value = 0.01895146864927834
$ws.Range("C4").Value = 0.8260900309969135
$ws.Range("D4").Value = 0.07508708545604498
$ws.Range("C5").Value = 0.555207060096766
$ws.Range("D5").Value = 0.0891774393803709
$ws.Range("C6").Value = 0.974622
$ws.Range("D6").Value = 0.01313232947922875
$ws.Range("C7").Value = 0.6600560089027289
$ws.Range("D7").Value = 0.074730904644722
$ws.Range("C8").Value = 0.892833505743799
$ws.Range("D8").Value = 0.02189482479449176
$ws.Range("C9").Value = 0.8011191828245071
$ws.Range("D9").Value = 0.04248145606026136
$ws.Range("C10").Value = 0.764908035263326
$ws.Range("D10").Value = 0.04486415568021546
$ws.Range("C11").Value = 0.6238796089843015
$ws.Range("D11").Value = 0.07859687146754231
$ws.Range("C12").Value = 0.91205
$ws.Range("D12").Value = 0.01640868672380577
$ws.Range("C13").Value = 0.764908035263326
$ws.Range("D13").Value = 0.04486415568021545

$ws = $wb.Worksheets.Item("KNN")
$ws.Range("C2").Value = 0.7163453738637385
$ws.Range("D2").Value = 0.04466161399986224
$ws.Range("C3").Value = 0.9075392670157069
$ws.Range("D3").Value = 0.01949278198374421
$ws.Range("C4").Value = 0.8079197795901624
$ws.Range("D4").Value = 0.07156676547607889
$ws.Range("C5").Value = 0.622206550802139
$ws.Range("D5").Value = 0.09835036864500163
$ws.Range("C6").Value = 0.9680319999999999
$ws.Range("D6").Value = 0.01398962355987246
$ws.Range("C7").Value = 0.6985813812129408
$ws.Range("D7").Value = 0.07494520892973978
$ws.Range("C8").Value = 0.9021972592452104
$ws.Range("D8").Value = 0.02202423046260341
$ws.Range("C9").Value = 0.8219568583248719
$ws.Range("D9").Value = 0.04280339391357577
$ws.Range("C10").Value = 0.7951162082125165
$ws.Range("D10").Value = 0.04887804631921126
$ws.Range("C11").Value = 0.6557439089796517
$ws.Range("D11").Value = 0.0803677482849833
$ws.Range("C12").Value = 0.9239599999999999
$ws.Range("D12").Value = 0.01832716284532835
$ws.Range("C13").Value = 0.7951162082125165
$ws.Range("D13").Value = 0.04887804631921126

$ws = $wb.Worksheets.Item("SVM")
$ws.Range("C2").Value = 0.7379692603421592
$ws.Range("D2").Value = 0.04292894322339021
$ws.Range("C3").Value = 0.8957068062827225
$ws.Range("D3").Value = 0.01678595080222053
$ws.Range("C4").Value = 0.8370894390060765
$ws.Range("D4").Value = 0.06694894545263737
$ws.Range("C5").Value = 0.5035385472370766
$ws.Range("D5").Value = 0.0890996442829423
$ws.Range("C6").Value = 0.9788020000000001
$ws.Range("D6").Value = 0.01065866496215684
$ws.Range("C7").Value = 0.624209215720661
$ws.Range("D7").Value = 0.07758895745878255
$ws.Range("C8").Value = 0.884279330367965
$ws.Range("D8").Value = 0.02108735563261568
$ws.Range("C9").Value = 0.7817992065139958
$ws.Range("D9").Value = 0.04316843797154366
$ws.Range("C10").Value = 0.7411684588040441
$ws.Range("D10").Value = 0.04399547776176307
$ws.Range("C11").Value = 0.5957184264445482
$ws.Range("D11").Value = 0.07508912676653019
$ws.Range("C12").Value = 0.9032119999999999
$ws.Range("D12").Value = 0.01568715372005815
$ws.Range("C13").Value = 0.7411684588040441
$ws.Range("D13").Value = 0.04399547776176306
